$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 898.6
$ws.Range("I6").Value = 190.6923
$ws.Range("K6").Value = 572.0769
$ws.Range("M6").Value = -460.0769
$ws.Range("H19").Value = 3728.9443
$ws.Range("I19").Value = 3699.6667
$ws.Range("K19").Value = 3699.6667
$ws.Range("M19").Value = -3524.6667
$ws.Range("H33").Value = 225.48276
$ws.Range("I33").Value = 219.76923
$ws.Range("J33").Value = 275
$ws.Range("K33").Value = 219.76923
$ws.Range("L33").Value = 275
$ws.Range("M33").Value = 9.230770000000007
$ws.Range("N33").Value = -733
$ws.Range("H96").Value = 45455796
$ws.Range("I96").Value = 806.26666
$ws.Range("K96").Value = 2418.79998
$ws.Range("M96").Value = -1045.79998
$ws.Range("H106").Value = 2691.8975
$ws.Range("I106").Value = 1805.5938
$ws.Range("K106").Value = 1805.5938
$ws.Range("M106").Value = -1174.5938
$ws.Range("H113").Value = 9898
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H116").Value = 5040.304
$ws.Range("J116").Value = 3831.6667
$ws.Range("L116").Value = 3831.6667
$ws.Range("N116").Value = -10715.6667
$ws.Range("H137").Value = 2251.3586
$ws.Range("I137").Value = 1844.4839
$ws.Range("J137").Value = 3092.2334
$ws.Range("K137").Value = 5533.4517
$ws.Range("L137").Value = 9276.700199999999
$ws.Range("M137").Value = -2983.4517
$ws.Range("N137").Value = -14376.7002
$ws.Range("H138").Value = 16134172
$ws.Range("J138").Value = 5253.316
$ws.Range("L138").Value = 15759.948
$ws.Range("N138").Value = -26039.948

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1917.4783
$ws.Range("I2").Value = 1909.6666
$ws.Range("K2").Value = 1909.6666
$ws.Range("M2").Value = -1796.6666
$ws.Range("H32").Value = 3394.5688
$ws.Range("I32").Value = 2723.4565
$ws.Range("J32").Value = 5967.1665
$ws.Range("K32").Value = 2723.4565
$ws.Range("L32").Value = 5967.1665
$ws.Range("M32").Value = -2436.4565
$ws.Range("N32").Value = -6541.1665
$ws.Range("H61").Value = 2116.75
$ws.Range("I61").Value = 1886.1714
$ws.Range("J61").Value = 3013.4443
$ws.Range("K61").Value = 1886.1714
$ws.Range("L61").Value = 3013.4443
$ws.Range("M61").Value = -1674.1714
$ws.Range("N61").Value = -3437.4443
$ws.Range("H97").Value = 3410.0417
$ws.Range("I97").Value = 4004.2632
$ws.Range("K97").Value = 4004.2632
$ws.Range("M97").Value = -3508.2632
$ws.Range("H110").Value = 1417.6471
$ws.Range("I110").Value = 1313.7667
$ws.Range("J110").Value = 2196.75
$ws.Range("K110").Value = 1313.7667
$ws.Range("L110").Value = 2196.75
$ws.Range("M110").Value = 731.2333000000001
$ws.Range("N110").Value = -6286.75
$ws.Range("H116").Value = 1917.4783
$ws.Range("I116").Value = 1909.6666
$ws.Range("K116").Value = 1909.6666
$ws.Range("M116").Value = 384.3334
$ws.Range("H122").Value = 6025.191
$ws.Range("I122").Value = 3842.1943
$ws.Range("K122").Value = 11526.5829
$ws.Range("M122").Value = -9076.582900000001
$ws.Range("H132").Value = 5279.639
$ws.Range("I132").Value = 4724
$ws.Range("K132").Value = 14172
$ws.Range("M132").Value = -11642
$ws.Range("H136").Value = 2116.75
$ws.Range("I136").Value = 1886.1714
$ws.Range("J136").Value = 3013.4443
$ws.Range("K136").Value = 5658.5142
$ws.Range("L136").Value = 9040.332900000001
$ws.Range("M136").Value = -3108.5142
$ws.Range("N136").Value = -14140.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1917.4783
$ws.Range("I3").Value = 1909.6666
$ws.Range("K3").Value = 1909.6666
$ws.Range("M3").Value = -1795.6666
$ws.Range("H94").Value = 2856.6287
$ws.Range("I94").Value = 1231.16
$ws.Range("K94").Value = 1231.16
$ws.Range("M94").Value = -780.1600000000001
$ws.Range("H99").Value = 10619.429
$ws.Range("J99").Value = 8800
$ws.Range("L99").Value = 8800
$ws.Range("N99").Value = -11796
$ws.Range("H134").Value = 6857.7
$ws.Range("I134").Value = 4511.2144
$ws.Range("K134").Value = 13533.6432
$ws.Range("M134").Value = -10998.6432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2687.3777
$ws.Range("I31").Value = 2631.7917
$ws.Range("J31").Value = 2750.9048
$ws.Range("K31").Value = 2631.7917
$ws.Range("L31").Value = 2750.9048
$ws.Range("M31").Value = -2336.7917
$ws.Range("N31").Value = -3340.9048
$ws.Range("H34").Value = 2687.3777
$ws.Range("I34").Value = 2631.7917
$ws.Range("J34").Value = 2750.9048
$ws.Range("K34").Value = 2631.7917
$ws.Range("L34").Value = 2750.9048
$ws.Range("M34").Value = -2429.7917
$ws.Range("N34").Value = -3154.9048
$ws.Range("H62").Value = 76928510
$ws.Range("I62").Value = 125004800
$ws.Range("J62").Value = 6460
$ws.Range("K62").Value = 125004800
$ws.Range("L62").Value = 6460
$ws.Range("M62").Value = -125004176
$ws.Range("N62").Value = -7708
$ws.Range("H65").Value = 76928510
$ws.Range("I65").Value = 125004800
$ws.Range("J65").Value = 6460
$ws.Range("K65").Value = 625024000
$ws.Range("L65").Value = 32300
$ws.Range("M65").Value = -625020880
$ws.Range("N65").Value = -38540
$ws.Range("H86").Value = 3709503.5
$ws.Range("I86").Value = 6066497.5
$ws.Range("K86").Value = 6066497.5
$ws.Range("M86").Value = -6065374.5
$ws.Range("H89").Value = 3709503.5
$ws.Range("I89").Value = 6066497.5
$ws.Range("K89").Value = 30332487.5
$ws.Range("M89").Value = -30326871.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 2383
$ws.Range("I86").Value = 545.375
$ws.Range("J86").Value = 4833.1665
$ws.Range("K86").Value = 1636.125
$ws.Range("L86").Value = 14499.4995
$ws.Range("M86").Value = -450.125
$ws.Range("N86").Value = -16871.4995
$ws.Range("H89").Value = 2383
$ws.Range("I89").Value = 545.375
$ws.Range("J89").Value = 4833.1665
$ws.Range("K89").Value = 4908.375
$ws.Range("L89").Value = 43498.4985
$ws.Range("M89").Value = 1019.625
$ws.Range("N89").Value = -55354.4985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 1308.3334
$ws.Range("I13").Value = 1645.1111
$ws.Range("K13").Value = 1645.1111
$ws.Range("M13").Value = -1506.1111
$ws.Range("H97").Value = 715.9
$ws.Range("I97").Value = 622.8
$ws.Range("K97").Value = 622.8
$ws.Range("M97").Value = -126.8
$ws.Range("H122").Value = 2442
$ws.Range("I122").Value = 2087.25
$ws.Range("K122").Value = 6261.75
$ws.Range("M122").Value = -3811.75
$ws.Range("H126").Value = 6480.9653
$ws.Range("I126").Value = 5011.3335
$ws.Range("K126").Value = 15034.0005
$ws.Range("M126").Value = -12564.0005
$ws.Range("H132").Value = 6730.3335
$ws.Range("I132").Value = 6525.143
$ws.Range("J132").Value = 7140.7144
$ws.Range("K132").Value = 19575.429
$ws.Range("L132").Value = 21422.1432
$ws.Range("M132").Value = -17045.429
$ws.Range("N132").Value = -26482.1432
$ws.Range("H134").Value = 267108.66
$ws.Range("J134").Value = 267108.66
$ws.Range("L134").Value = 801325.98
$ws.Range("N134").Value = -806395.98

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2069.125
$ws.Range("J93").Value = 914
$ws.Range("L93").Value = 914
$ws.Range("N93").Value = -3410
$ws.Range("H100").Value = 2831.818
$ws.Range("I100").Value = 1700
$ws.Range("K100").Value = 1700
$ws.Range("M100").Value = -1159
$ws.Range("H122").Value = 4940.8965
$ws.Range("I122").Value = 3619.95
$ws.Range("K122").Value = 10859.85
$ws.Range("M122").Value = -8409.849999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3055.875
$ws.Range("I96").Value = 3278.2856
$ws.Range("K96").Value = 3278.2856
$ws.Range("M96").Value = -1905.2856
$ws.Range("H100").Value = 720.9091
$ws.Range("I100").Value = 605.17645
$ws.Range("J100").Value = 1114.4
$ws.Range("K100").Value = 1210.3529
$ws.Range("L100").Value = 2228.8
$ws.Range("M100").Value = -669.3529000000001
$ws.Range("N100").Value = -3310.8
$ws.Range("H107").Value = 1799.8334
$ws.Range("I107").Value = 1745.1666
$ws.Range("J107").Value = 1881.8334
$ws.Range("K107").Value = 5235.4998
$ws.Range("L107").Value = 5645.5002
$ws.Range("M107").Value = -3315.4998
$ws.Range("N107").Value = -9485.5002
$ws.Range("H139").Value = 69943.14
$ws.Range("J139").Value = 69943.14
$ws.Range("L139").Value = 69943.14
$ws.Range("N139").Value = -80223.14
